# Actualización automática hashcode dom sep 20 01:29:39 CEST 2020
#
# This script replicates the "hashcode" column refresh performed by the
# original commit: a set of previously computed hash values (column B)
# are recomputed/updated while the identifiers in column A stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B34").Value  = "1a2aad99247432a7c8ad2c855eaeec1e"
$ws.Range("B154").Value = "6b15316edc1cc092b4abac42be90bd28"
$ws.Range("B160").Value = "a971ea9eb8c3823f3586968e3793190b"
$ws.Range("B162").Value = "b2958ca0a2f48c38ed413b0942283382"
$ws.Range("B169").Value = "4da83de0fa8baa0c3e34ef948fa497bf"
$ws.Range("B180").Value = "9ff250cc2296e8b04e2e9c55eb7b492a"
$ws.Range("B213").Value = "289d9c7f686850f0271f00b042591a5a"
$ws.Range("B222").Value = "d0a510c33c0ac6bb6a7521f08fec4070"
$ws.Range("B227").Value = "82760c335d1800fd1aeb50687d6f826e"
$ws.Range("B229").Value = "9a8cc75de1629534c3eaece5b8c32057"
$ws.Range("B232").Value = "3f0a589ba5292d038af5d7e15f995d2b"
$ws.Range("B284").Value = "afc91a4d0896544a39504d970bebe301"
$ws.Range("B468").Value = "76fe75e6b689c434da60d249ba6765bf"
$ws.Range("B486").Value = "090ce60a84e4df080ad7c313bf00d29a"
$ws.Range("B516").Value = "3573f972709eca56275fd504bb286c75"
$ws.Range("B524").Value = "e3d6f2571a6e47a237de56acc60583d0"
$ws.Range("B535").Value = "c2ff6a83c1beba8689e2d6eaa3eb06e1"
$ws.Range("B545").Value = "caed40e30b8d326c9ee29159f49801d9"
$ws.Range("B565").Value = "6dae6fa19d878e3e786208dc34f13627"
$ws.Range("B578").Value = "c2773ef09b571a4d55e3f514b1138e7d"
$ws.Range("B678").Value = "7f37c26eae181fa0ad2e97b5864751b2"
$ws.Range("B692").Value = "4fc5fa4b3dd3ce2d2f863a4ac7f1255b"
$ws.Range("B697").Value = "536052429b70078e1e780ef554fbc516"
$ws.Range("B712").Value = "32cabfb6d54c47197f02bfa132f2bceb"
$ws.Range("B715").Value = "fb3404a2ee3af1938e8f92d2e045b730"
$ws.Range("B741").Value = "93049bfcc2ff1ccbc37fcd3a7fe75f92"
$ws.Range("B823").Value = "d05f60cb7fe7ed68b218c83ac767a514"
$ws.Range("B827").Value = "828dfcdbe017b46b27ba6a91372baea2"
